$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 5, shifting rows 6:11 up to become rows 5:10
$ws.Rows.Item(5).Delete()

# Update the active selection to B6, matching the post-edit worksheet view
$ws.Range("B6").Select()
